# Auto-generated Excel COM-interop script
# Applies the scheduled-runner's refreshed market-price / profit figures
# (currentAveragePrice*, LevePriceNQ/HQ, LeveProfitNQ/HQ) across all 8 job sheets,
# matching the upstream OOXML diff cell-for-cell.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 695.8125
$ws.Range("I2").Value = 309.3
$ws.Range("J2").Value = 1340
$ws.Range("K2").Value = 309.3
$ws.Range("L2").Value = 1340
$ws.Range("M2").Value = -196.3
$ws.Range("N2").Value = -1566
$ws.Range("H3").Value = 33333
$ws.Range("J3").Value = 33333
$ws.Range("L3").Value = 33333
$ws.Range("N3").Value = -33561
$ws.Range("H33").Value = 143467.58
$ws.Range("I33").Value = 205.7
$ws.Range("K33").Value = 205.7
$ws.Range("M33").Value = 23.30000000000001
$ws.Range("H80").Value = 835.9091
$ws.Range("J80").Value = 1070.7142
$ws.Range("L80").Value = 3212.1426
$ws.Range("N80").Value = -5208.142599999999
$ws.Range("H83").Value = 835.9091
$ws.Range("J83").Value = 1070.7142
$ws.Range("L83").Value = 9636.427799999999
$ws.Range("N83").Value = -19620.4278
$ws.Range("H98").Value = 2432.5334
$ws.Range("I98").Value = 732.5
$ws.Range("J98").Value = 9232.666999999999
$ws.Range("K98").Value = 732.5
$ws.Range("L98").Value = 9232.666999999999
$ws.Range("M98").Value = 765.5
$ws.Range("N98").Value = -12228.667
$ws.Range("H102").Value = 33333
$ws.Range("J102").Value = 33333
$ws.Range("L102").Value = 33333
$ws.Range("N102").Value = -39823
$ws.Range("H113").Value = 3272
$ws.Range("I113").Value = 2424.125
$ws.Range("K113").Value = 2424.125
$ws.Range("M113").Value = 829.875
$ws.Range("H122").Value = 2432.5334
$ws.Range("I122").Value = 732.5
$ws.Range("J122").Value = 9232.666999999999
$ws.Range("K122").Value = 2197.5
$ws.Range("L122").Value = 27698.001
$ws.Range("M122").Value = 252.5
$ws.Range("N122").Value = -32598.001
$ws.Range("H132").Value = 3078.8982
$ws.Range("I132").Value = 2721.018
$ws.Range("K132").Value = 8163.054
$ws.Range("M132").Value = -5633.054

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4538.1333
$ws.Range("I32").Value = 903.9184
$ws.Range("K32").Value = 903.9184
$ws.Range("M32").Value = -616.9184
$ws.Range("H61").Value = 4942.242
$ws.Range("I61").Value = 4022.1482
$ws.Range("J61").Value = 9082.666999999999
$ws.Range("K61").Value = 4022.1482
$ws.Range("L61").Value = 9082.666999999999
$ws.Range("M61").Value = -3810.1482
$ws.Range("N61").Value = -9506.666999999999
$ws.Range("H74").Value = 4776
$ws.Range("I74").Value = 3784.6667
$ws.Range("K74").Value = 3784.6667
$ws.Range("M74").Value = -2910.6667
$ws.Range("H77").Value = 4776
$ws.Range("I77").Value = 3784.6667
$ws.Range("K77").Value = 18923.3335
$ws.Range("M77").Value = -14555.3335
$ws.Range("H122").Value = 5197.028
$ws.Range("I122").Value = 4644.8623
$ws.Range("K122").Value = 13934.5869
$ws.Range("M122").Value = -11484.5869
$ws.Range("H136").Value = 4942.242
$ws.Range("I136").Value = 4022.1482
$ws.Range("J136").Value = 9082.666999999999
$ws.Range("K136").Value = 12066.4446
$ws.Range("L136").Value = 27248.001
$ws.Range("M136").Value = -9516.444600000001
$ws.Range("N136").Value = -32348.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4631244.5
$ws.Range("I86").Value = 6411177.5
$ws.Range("K86").Value = 6411177.5
$ws.Range("M86").Value = -6410054.5
$ws.Range("H89").Value = 4631244.5
$ws.Range("I89").Value = 6411177.5
$ws.Range("K89").Value = 32055887.5
$ws.Range("M89").Value = -32050271.5
$ws.Range("H94").Value = 30021052
$ws.Range("I94").Value = 41667544
$ws.Range("J94").Value = 72931
$ws.Range("K94").Value = 41667544
$ws.Range("L94").Value = 72931
$ws.Range("M94").Value = -41667093
$ws.Range("N94").Value = -73833
$ws.Range("H99").Value = 1883.4348
$ws.Range("I99").Value = 1918.5454
$ws.Range("K99").Value = 1918.5454
$ws.Range("M99").Value = -420.5454
$ws.Range("H107").Value = 2840.6924
$ws.Range("I107").Value = 2567.5715
$ws.Range("J107").Value = 3987.8
$ws.Range("K107").Value = 2567.5715
$ws.Range("L107").Value = 3987.8
$ws.Range("M107").Value = -647.5715
$ws.Range("N107").Value = -7827.8
$ws.Range("H132").Value = 64611.668
$ws.Range("J132").Value = 64611.668
$ws.Range("L132").Value = 64611.668
$ws.Range("N132").Value = -74731.66800000001
$ws.Range("H134").Value = 3593.0527
$ws.Range("I134").Value = 3593.0527
$ws.Range("K134").Value = 10779.1581
$ws.Range("M134").Value = -8244.158100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6616.8667
$ws.Range("I99").Value = 5470.875
$ws.Range("K99").Value = 5470.875
$ws.Range("M99").Value = -3972.875
$ws.Range("H126").Value = 6616.8667
$ws.Range("I126").Value = 5470.875
$ws.Range("K126").Value = 16412.625
$ws.Range("M126").Value = -13942.625
$ws.Range("H140").Value = 68780
$ws.Range("J140").Value = 68780
$ws.Range("L140").Value = 68780
$ws.Range("N140").Value = -79140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 8253.454
$ws.Range("I7").Value = 10085.714
$ws.Range("J7").Value = 5047
$ws.Range("K7").Value = 30257.142
$ws.Range("L7").Value = 15141
$ws.Range("M7").Value = -30145.142
$ws.Range("N7").Value = -15365
$ws.Range("H12").Value = 380.25
$ws.Range("I12").Value = 173
$ws.Range("K12").Value = 519
$ws.Range("M12").Value = -346
$ws.Range("H23").Value = 489.4
$ws.Range("I23").Value = 35
$ws.Range("J23").Value = 603
$ws.Range("K23").Value = 105
$ws.Range("L23").Value = 1809
$ws.Range("M23").Value = 130
$ws.Range("N23").Value = -2279
$ws.Range("H88").Value = 4302.8
$ws.Range("I88").Value = 1914
$ws.Range("K88").Value = 5742
$ws.Range("M88").Value = -5314
$ws.Range("H91").Value = 4302.8
$ws.Range("I91").Value = 1914
$ws.Range("K91").Value = 5742
$ws.Range("M91").Value = -4260
$ws.Range("H134").Value = 4329.5
$ws.Range("I134").Value = 1215.4667
$ws.Range("J134").Value = 19899.666
$ws.Range("K134").Value = 3646.4001
$ws.Range("L134").Value = 59698.99800000001
$ws.Range("M134").Value = 1423.5999
$ws.Range("N134").Value = -69838.99800000001
$ws.Range("H137").Value = 3747.2856
$ws.Range("J137").Value = 3756.1
$ws.Range("L137").Value = 11268.3
$ws.Range("N137").Value = -21468.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 244.6
$ws.Range("I2").Value = 244.6
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 244.6
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -131.6
$ws.Range("N2").ClearContents()
$ws.Range("H102").Value = 1763.4688
$ws.Range("I102").Value = 1312.2963
$ws.Range("K102").Value = 1312.2963
$ws.Range("M102").Value = 309.7037
$ws.Range("H122").Value = 5732.85
$ws.Range("I122").Value = 5400.5454
$ws.Range("J122").Value = 6139
$ws.Range("K122").Value = 16201.6362
$ws.Range("L122").Value = 18417
$ws.Range("M122").Value = -13751.6362
$ws.Range("N122").Value = -23317
$ws.Range("H126").Value = 6866.8887
$ws.Range("I126").Value = 5778.6665
$ws.Range("K126").Value = 17335.9995
$ws.Range("M126").Value = -14865.9995
$ws.Range("H132").Value = 4881.892
$ws.Range("I132").Value = 3842.1177
$ws.Range("K132").Value = 11526.3531
$ws.Range("M132").Value = -8996.3531
$ws.Range("H139").Value = 81858
$ws.Range("J139").Value = 81858
$ws.Range("L139").Value = 81858
$ws.Range("N139").Value = -92138

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1273.5
$ws.Range("I61").Value = 1567.25
$ws.Range("K61").Value = 1567.25
$ws.Range("M61").Value = -1365.25
$ws.Range("H100").Value = 1950
$ws.Range("I100").Value = 1975
$ws.Range("K100").Value = 1975
$ws.Range("M100").Value = -1434
$ws.Range("H113").Value = 1273.5
$ws.Range("I113").Value = 1567.25
$ws.Range("K113").Value = 1567.25
$ws.Range("M113").Value = 602.75
$ws.Range("H122").Value = 3499.625
$ws.Range("I122").Value = 2499.25
$ws.Range("K122").Value = 7497.75
$ws.Range("M122").Value = -5047.75
$ws.Range("H132").Value = 2397.7273
$ws.Range("I132").Value = 2397.7273
$ws.Range("K132").Value = 7193.1819
$ws.Range("M132").Value = -4663.1819
$ws.Range("H140").Value = 133985.75
$ws.Range("J140").Value = 133664.33
$ws.Range("L140").Value = 133664.33
$ws.Range("N140").Value = -144024.33

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 320548.3
$ws.Range("J14").Value = 2476749.5
$ws.Range("L14").Value = 2476749.5
$ws.Range("N14").Value = -2477085.5
$ws.Range("H119").Value = 1529924.2
$ws.Range("J119").Value = 39899
$ws.Range("L119").Value = 39899
$ws.Range("N119").Value = -49575
$ws.Range("H122").Value = 3331.8948
$ws.Range("I122").Value = 3214.8667
$ws.Range("K122").Value = 9644.6001
$ws.Range("M122").Value = -7194.6001
$ws.Range("H126").Value = 4344.4546
$ws.Range("I126").Value = 2348.3333
$ws.Range("K126").Value = 7044.999899999999
$ws.Range("M126").Value = -4574.999899999999
$ws.Range("H132").Value = 5922.8
$ws.Range("I132").Value = 3911.9143
$ws.Range("J132").Value = 19999
$ws.Range("K132").Value = 11735.7429
$ws.Range("L132").Value = 59997
$ws.Range("M132").Value = -9205.742899999999
$ws.Range("N132").Value = -65057
$ws.Range("H136").Value = 5290.517
$ws.Range("I136").Value = 5297.174
$ws.Range("J136").Value = 5265
$ws.Range("K136").Value = 15891.522
$ws.Range("L136").Value = 15795
$ws.Range("M136").Value = -13341.522
$ws.Range("N136").Value = -20895
